$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new skill row (id 18 / NaturalLens / PassiveSkill / 0 / 0)
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "NaturalLens"
$ws.Range("C20").Value = "PassiveSkill"
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0

# Move the active selection like the author's click to C21 after entering the row
$ws.Range("C21").Select()
